# Update "想去人数" (want-to-go count) values in column F across the four
# worksheets (展览, 演出, 本地生活, 全部类型) to match the latest generated output.

$wb = $excel.ActiveWorkbook

$updates = @(
    @{Sheet = "展览";   Row = 2;  Value = 51},
    @{Sheet = "展览";   Row = 3;  Value = 259},
    @{Sheet = "展览";   Row = 4;  Value = 855},
    @{Sheet = "展览";   Row = 6;  Value = 425},
    @{Sheet = "展览";   Row = 7;  Value = 626},
    @{Sheet = "展览";   Row = 8;  Value = 229},
    @{Sheet = "展览";   Row = 11; Value = 164},
    @{Sheet = "展览";   Row = 12; Value = 732},
    @{Sheet = "展览";   Row = 13; Value = 98},
    @{Sheet = "展览";   Row = 14; Value = 1854},
    @{Sheet = "展览";   Row = 15; Value = 384},
    @{Sheet = "展览";   Row = 16; Value = 4247},
    @{Sheet = "展览";   Row = 17; Value = 386},
    @{Sheet = "展览";   Row = 19; Value = 14},
    @{Sheet = "展览";   Row = 20; Value = 64},
    @{Sheet = "展览";   Row = 21; Value = 152},
    @{Sheet = "展览";   Row = 22; Value = 132},

    @{Sheet = "演出";   Row = 4;  Value = 46},
    @{Sheet = "演出";   Row = 7;  Value = 484},

    @{Sheet = "本地生活"; Row = 2;  Value = 5387},
    @{Sheet = "本地生活"; Row = 3;  Value = 336},
    @{Sheet = "本地生活"; Row = 4;  Value = 303},

    @{Sheet = "全部类型"; Row = 2;  Value = 51},
    @{Sheet = "全部类型"; Row = 3;  Value = 5387},
    @{Sheet = "全部类型"; Row = 4;  Value = 336},
    @{Sheet = "全部类型"; Row = 6;  Value = 303},
    @{Sheet = "全部类型"; Row = 7;  Value = 259},
    @{Sheet = "全部类型"; Row = 9;  Value = 46},
    @{Sheet = "全部类型"; Row = 12; Value = 484},
    @{Sheet = "全部类型"; Row = 13; Value = 855},
    @{Sheet = "全部类型"; Row = 17; Value = 425},
    @{Sheet = "全部类型"; Row = 18; Value = 626},
    @{Sheet = "全部类型"; Row = 19; Value = 229},
    @{Sheet = "全部类型"; Row = 23; Value = 164},
    @{Sheet = "全部类型"; Row = 26; Value = 732},
    @{Sheet = "全部类型"; Row = 27; Value = 98},
    @{Sheet = "全部类型"; Row = 29; Value = 1854},
    @{Sheet = "全部类型"; Row = 30; Value = 384},
    @{Sheet = "全部类型"; Row = 31; Value = 4247},
    @{Sheet = "全部类型"; Row = 33; Value = 386},
    @{Sheet = "全部类型"; Row = 35; Value = 14},
    @{Sheet = "全部类型"; Row = 36; Value = 64},
    @{Sheet = "全部类型"; Row = 38; Value = 152},
    @{Sheet = "全部类型"; Row = 40; Value = 132}
)

foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item($u.Sheet)
    $ws.Cells.Item($u.Row, 6).Value = $u.Value
}
